$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Tabelle1" is already the active/selected sheet

# Week 13 ("28") had no actual hours recorded yet; fill in the real value.
# The Diff column (G13 = F13 - E13, a shared formula) recalculates on its own.
$ws.Range("F13").Value = 16

# The "Diff" total (G25) was never filled in - add the column sum, mirroring
# the existing totals in E25/F25.
$ws.Range("G25").Formula = "=SUM(G2:G24)"

# Make sure every cached formula result (G13, F25, G25, F26, ...) is fresh.
$excel.CalculateFull()

# Leave the selection on the cell that was just edited.
$ws.Range("G25").Select()
